$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.840.67"
$ws.Range("E2").Value = "  -1.87%  "
$ws.Range("D3").Value = "1.800.89"
$ws.Range("E3").Value = "  -1.36%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'309.14"
$ws.Range("E5").Value = "  -1.77%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").Value = "'0.4669"
$ws.Range("E7").Value = "  +4.26%  "
$ws.Range("D8").Value = "'0.3696"
$ws.Range("E8").Value = "  -2.02%  "
$ws.Range("D9").Value = "'0.07376"
$ws.Range("E9").Value = "  -0.98%  "
$ws.Range("D10").Value = "'0.8700"
$ws.Range("E10").Value = "  -2.06%  "
$ws.Range("E11").Value = "  -3.05%  "
$ws.Range("D12").Value = "1.855.49"
$ws.Range("E12").Value = "  +1.60%  "
$ws.Range("D13").Value = "'5.355"
$ws.Range("E13").Value = "  -1.97%  "
$ws.Range("D14").Value = "'92.28"
$ws.Range("E14").Value = "  -1.45%  "
$ws.Range("D15").Value = "'6.487"
$ws.Range("E15").Value = "  -3.79%  "
$ws.Range("D16").Value = "'0.07033"
$ws.Range("E16").Value = "  -1.49%  "
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("D18").Value = "'0.000008694"
$ws.Range("E18").Value = "  -1.17%  "
$ws.Range("E19").Value = "  -0.05%  "
$ws.Range("E20").Value = "  -3.15%  "
$ws.Range("D21").Value = "26.837.82"
$ws.Range("E21").Value = "  -1.91%  "
$ws.Range("D22").Value = "'5.291"
$ws.Range("E22").Value = "  -2.03%  "
$ws.Range("D23").Value = "'10.61"
$ws.Range("E23").Value = "  -3.38%  "
$ws.Range("D24").Value = "2.001.44"
$ws.Range("E24").Value = "  -2.50%  "
$ws.Range("D25").Value = "'1.891"
$ws.Range("E25").Value = "  -4.11%  "
$ws.Range("D26").Value = "'151.61"
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").Value = "'18.35"
$ws.Range("E27").Value = "  -1.91%  "
$ws.Range("D28").Value = "'2.138"
$ws.Range("E28").Value = "  -8.12%  "
$ws.Range("D29").Value = "'5.259"
$ws.Range("E29").Value = "  -2.39%  "
$ws.Range("D30").Value = "'115.91"
$ws.Range("D31").Value = "'0.08931"
$ws.Range("E31").Value = "  +0.44%  "
$ws.Range("D32").Value = "'0.7585"
$ws.Range("E32").Value = "  -4.11%  "
$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D33").Value = "'2.928"
$ws.Range("E33").Value = "  +0.29%  "
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "'1.148"
$ws.Range("E34").Value = "  -4.44%  "
$ws.Range("D35").Value = "'4.457"
$ws.Range("E35").Value = "  -3.37%  "
$ws.Range("E37").Value = "  -1.03%  "
$ws.Range("D38").Value = "'0.01957"
$ws.Range("E38").Value = "  -1.58%  "
$ws.Range("D39").Value = "'0.05248"
$ws.Range("E39").Value = "  -1.15%  "
$ws.Range("D40").Value = "'2.928"
$ws.Range("E40").Value = "  +1.78%  "
$ws.Range("D41").Value = "'7.223"
$ws.Range("E41").Value = "  -1.09%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "'2.378"
$ws.Range("E42").Value = "  +2.52%  "
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "'0.5290"
$ws.Range("E43").Value = "  -1.09%  "
$ws.Range("D44").Value = "'0.1660"
$ws.Range("E44").Value = "  -3.46%  "
$ws.Range("D45").Value = "'8.493"
$ws.Range("E45").Value = "  -2.04%  "
$ws.Range("D46").Value = "'0.5003"
$ws.Range("E46").Value = "  -1.95%  "
$ws.Range("D47").Value = "'10.26"
$ws.Range("E47").Value = "  -3.72%  "
$ws.Range("D48").Value = "'104.12"
$ws.Range("E48").Value = "  -1.05%  "
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("D50").Value = "'1.663"
$ws.Range("E50").Value = "  -2.16%  "
$ws.Range("D51").Value = "'0.06284"
$ws.Range("E51").Value = "  -1.91%  "
